$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

$ws_ALC.Range("H11").Value = 2705.3333
$ws_ALC.Range("I11").Value = 2705.3333
$ws_ALC.Range("K11").Value = 2705.3333
$ws_ALC.Range("M11").Value = -2565.3333

$ws_ALC.Range("H61").Value = 3175
$ws_ALC.Range("I61").Value = 5300
$ws_ALC.Range("J61").Value = 1050
$ws_ALC.Range("K61").Value = 15900
$ws_ALC.Range("L61").Value = 3150
$ws_ALC.Range("M61").Value = -15728
$ws_ALC.Range("N61").Value = -3494

$ws_ALC.Range("H107").Value = 28754084
$ws_ALC.Range("J107").Value = 35003170
$ws_ALC.Range("L107").Value = 35003170
$ws_ALC.Range("N107").Value = -35007010

$ws_ALC.Range("H132").Value = 2269.4614
$ws_ALC.Range("I132").Value = 2226.96
$ws_ALC.Range("K132").Value = 6680.88
$ws_ALC.Range("M132").Value = -4150.88

$ws_ALC.Range("H137").Value = 3936.8064
$ws_ALC.Range("I137").Value = 4243.6816
$ws_ALC.Range("J137").Value = 3186.6667
$ws_ALC.Range("K137").Value = 12731.0448
$ws_ALC.Range("L137").Value = 9560.000100000001
$ws_ALC.Range("M137").Value = -10181.0448
$ws_ALC.Range("N137").Value = -14660.0001

$ws_ALC.Range("H141").Value = 3410.5293
$ws_ALC.Range("I141").Value = 3260.862
$ws_ALC.Range("K141").Value = 9782.585999999999
$ws_ALC.Range("M141").Value = -4602.585999999999

$ws_ARM.Range("H32").Value = 1456254.5
$ws_ARM.Range("I32").Value = 1669445.6
$ws_ARM.Range("K32").Value = 1669445.6
$ws_ARM.Range("M32").Value = -1669158.6

$ws_ARM.Range("H45").Value = 3270.2856
$ws_ARM.Range("I45").Value = 2866
$ws_ARM.Range("K45").Value = 2866
$ws_ARM.Range("M45").Value = -2489

$ws_ARM.Range("H61").Value = 6380.3335
$ws_ARM.Range("I61").Value = 3183.1875
$ws_ARM.Range("J61").Value = 12774.625
$ws_ARM.Range("K61").Value = 3183.1875
$ws_ARM.Range("L61").Value = 12774.625
$ws_ARM.Range("M61").Value = -2971.1875
$ws_ARM.Range("N61").Value = -13198.625

$ws_ARM.Range("H74").Value = 61205.895
$ws_ARM.Range("I74").Value = 90875.836
$ws_ARM.Range("K74").Value = 90875.836
$ws_ARM.Range("M74").Value = -90001.836

$ws_ARM.Range("H77").Value = 61205.895
$ws_ARM.Range("I77").Value = 90875.836
$ws_ARM.Range("K77").Value = 454379.18
$ws_ARM.Range("M77").Value = -450011.18

$ws_ARM.Range("H132").Value = 1653116.6
$ws_ARM.Range("I132").Value = 4394492
$ws_ARM.Range("J132").Value = 8291.35
$ws_ARM.Range("K132").Value = 13183476
$ws_ARM.Range("L132").Value = 24874.05
$ws_ARM.Range("M132").Value = -13180946
$ws_ARM.Range("N132").Value = -29934.05

$ws_ARM.Range("H136").Value = 6380.3335
$ws_ARM.Range("I136").Value = 3183.1875
$ws_ARM.Range("J136").Value = 12774.625
$ws_ARM.Range("K136").Value = 9549.5625
$ws_ARM.Range("L136").Value = 38323.875
$ws_ARM.Range("M136").Value = -6999.5625
$ws_ARM.Range("N136").Value = -43423.875

$ws_BSM.Range("H22").Value = 4831.773
$ws_BSM.Range("I22").Value = 6589.375
$ws_BSM.Range("K22").Value = 6589.375
$ws_BSM.Range("M22").Value = -6416.375

$ws_BSM.Range("H99").Value = 4548075
$ws_BSM.Range("I99").Value = 3164.0667
$ws_BSM.Range("K99").Value = 3164.0667
$ws_BSM.Range("M99").Value = -1666.0667

$ws_CRP.Range("H31").Value = 8981.375
$ws_CRP.Range("I31").Value = 3758
$ws_CRP.Range("K31").Value = 3758
$ws_CRP.Range("M31").Value = -3463

$ws_CRP.Range("H34").Value = 8981.375
$ws_CRP.Range("I34").Value = 3758
$ws_CRP.Range("K34").Value = 3758
$ws_CRP.Range("M34").Value = -3556

$ws_CRP.Range("H132").Value = 6810.263
$ws_CRP.Range("I132").Value = 2240
$ws_CRP.Range("K132").Value = 6720
$ws_CRP.Range("M132").Value = -4190

$ws_CRP.Range("H134").Value = 10869.692
$ws_CRP.Range("I134").Value = 4202.75
$ws_CRP.Range("K134").Value = 12608.25
$ws_CRP.Range("M134").Value = -10073.25

$ws_CUL.Range("H2").Value = 201410.6
$ws_CUL.Range("I2").Value = 484
$ws_CUL.Range("K2").Value = 2904
$ws_CUL.Range("M2").Value = -2791

$ws_CUL.Range("H55").Value = 43340384
$ws_CUL.Range("J55").Value = 25016248
$ws_CUL.Range("L55").Value = 75048744
$ws_CUL.Range("N55").Value = -75049098

$ws_CUL.Range("H121").Value = 4630.4287
$ws_CUL.Range("J121").Value = 4630.4287
$ws_CUL.Range("L121").Value = 13891.2861
$ws_CUL.Range("N121").Value = -16511.2861

$ws_CUL.Range("H134").Value = 232639.81
$ws_CUL.Range("I134").Value = 232639.81
$ws_CUL.Range("K134").Value = 697919.4299999999
$ws_CUL.Range("M134").Value = -692849.4299999999

$ws_GSM.Range("H80").Value = 102299.2
$ws_GSM.Range("I80").Value = 2083
$ws_GSM.Range("K80").Value = 2083
$ws_GSM.Range("M80").Value = -1085

$ws_GSM.Range("H83").Value = 102299.2
$ws_GSM.Range("I83").Value = 2083
$ws_GSM.Range("K83").Value = 10415
$ws_GSM.Range("M83").Value = -5423

$ws_GSM.Range("H102").Value = 1996.7273
$ws_GSM.Range("I102").Value = 1594.5555
$ws_GSM.Range("K102").Value = 1594.5555
$ws_GSM.Range("M102").Value = 27.44450000000006

$ws_GSM.Range("H113").Value = 6175.2905
$ws_GSM.Range("I113").Value = 2695.6
$ws_GSM.Range("K113").Value = 2695.6
$ws_GSM.Range("M113").Value = -525.5999999999999

$ws_GSM.Range("H132").Value = 8267
$ws_GSM.Range("I132").Value = 3832.3333
$ws_GSM.Range("J132").Value = 13588.6
$ws_GSM.Range("K132").Value = 11496.9999
$ws_GSM.Range("L132").Value = 40765.8
$ws_GSM.Range("M132").Value = -8966.999899999999
$ws_GSM.Range("N132").Value = -45825.8

$ws_LTW.Range("H93").Value = 3960.125
$ws_LTW.Range("I93").Value = 3353.4092
$ws_LTW.Range("K93").Value = 3353.4092
$ws_LTW.Range("M93").Value = -2105.4092

$ws_LTW.Range("H100").Value = 4659.364
$ws_LTW.Range("I100").Value = 4138.778
$ws_LTW.Range("K100").Value = 4138.778
$ws_LTW.Range("M100").Value = -3597.778

$ws_LTW.Range("H132").Value = 16136368
$ws_LTW.Range("I132").Value = 33336464
$ws_LTW.Range("K132").Value = 100009392
$ws_LTW.Range("M132").Value = -100006862

$ws_LTW.Range("H136").Value = 10808.419
$ws_LTW.Range("I136").Value = 7614.5557
$ws_LTW.Range("J136").Value = 15230.692
$ws_LTW.Range("K136").Value = 22843.6671
$ws_LTW.Range("L136").Value = 45692.076
$ws_LTW.Range("M136").Value = -20293.6671
$ws_LTW.Range("N136").Value = -50792.076

$ws_WVR.Range("H132").Value = 20020616
$ws_WVR.Range("I132").Value = 22738564
$ws_WVR.Range("K132").Value = 68215692
$ws_WVR.Range("M132").Value = -68213162

$ws_WVR.Range("H136").Value = 40045570
$ws_WVR.Range("I136").Value = 76925630
$ws_WVR.Range("J136").Value = 92166.586
$ws_WVR.Range("K136").Value = 230776890
$ws_WVR.Range("L136").Value = 276499.758
$ws_WVR.Range("M136").Value = -230774340
$ws_WVR.Range("N136").Value = -281599.758
